$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ofsted_cs_send_inspections")

# Column K = inspection_theme. Update rows 2-60 to match corrected capitalisation
# (sentence case: first word capitalised, rest lower-cased; a few rows also had
# trailing narrative text trimmed off).
$updates = @{
    2 = 'Identification of initial need and risk in barnsley.'
    3 = 'Identification of initial need and risk in bedford borough council for children and families who need help.'
    4 = 'Identification of initial need and risk in blackburn with darwen.'
    5 = 'Sexual abuse in the family in bracknell forest between 21 and 25 january 2019.'
    6 = 'Abuse and neglect in bristol city council between 16 and 20 october 2017.'
    7 = 'Identification of initial need and risk in buckinghamshire.'
    8 = 'Abuse and neglect in central bedfordshire between 14 and 18 march 2016.'
    9 = 'The criminal exploitation of children in cheshire east.'
    10 = 'Abuse and neglect in cheshire west and chester between 25 and 29 september 2017.'
    11 = 'Abuse and neglect in bradford metropolitan district council between 27 february and 3 march 2017.'
    12 = 'Sexual abuse in the family in york between 24 and 28 september 2018.'
    13 = 'Sexual abuse in the family in cornwall between 8 october 2018 and 12 october 2018.'
    14 = 'Serious youth violence in coventry.'
    15 = 'Abuse and neglect in derby city between 18 and 22 march 2019.'
    16 = 'Domestic abuse in durham between 9 and 13 july 2018.'
    17 = 'Abuse and neglect in east sussex.'
    18 = 'Identification of initial need and risk in gloucestershire.'
    19 = 'Child exploitation in halton between 8 and 12 july 2019.'
    20 = 'Abuse and neglect in hampshire between 5 and 9 december 2016.'
    21 = 'The criminal exploitation of children in kirklees.'
    22 = 'Serious youth violence in lancashire.'
    23 = 'Serious youth violence in leeds.'
    24 = 'Abuse and neglect in lincolnshire between 17 and 21 october 2016.'
    25 = 'Abuse and neglect in liverpool between 20 and 24 june.'
    26 = 'Childrens mental health in bexley between 20 and 24 january 2020.'
    27 = 'Abuse and neglect in the london borough of croydon between 16 may and 20 may.'
    28 = 'Child sexual exploitation, children associated with gangs and at risk of exploitation and children missing from home, care or education between 12 and 16 february 2018.'
    29 = 'Abuse and neglect in haringey between 4 and 8 december 2017.'
    30 = 'Children and families who need help in harrow.'
    31 = 'Abuse and neglect in hounslow between 21 march 2017 and 24 march 2017.'
    32 = 'Sexual abuse in the family in islington between 3 december 2018 and 7 december 2018.'
    33 = 'Initial need and risk in the london borough of lewisham.'
    34 = 'Serious youth violence in merton.'
    35 = 'Identification of initial need and risk (often referred to as the front door) in sutton.'
    36 = 'Serious youth violence in manchester.'
    37 = 'Domestic abuse in medway between 18 june 2018 and 22 june 2018.'
    38 = 'Childrens mental health in milton keynes between 14 and 18 october 2019.'
    39 = 'Child exploitation in northumberland between 17 and 21 june 2019.'
    40 = 'Abuse and neglect in oxfordshire between 7 march 2016 and 12 march 2016.'
    41 = 'Abuse and neglect in peterborough city council between 26 and 30 june 2017.'
    42 = 'Childrens mental health in plymouth between 18 november 2019 and 22 november 2019.'
    43 = 'Childrens mental health in portsmouth between 9 and 13 december 2019.'
    44 = 'Identification of initial need and risk (often referred to as the front door) in rochdale.'
    45 = 'Identification of initial need and risk in the royal borough of windsor and maidenhead.'
    46 = 'Abuse and neglect in salford between 12 and 16 september 2016.'
    47 = 'Childrens mental health in sefton.'
    48 = 'Sexual abuse in the family in shropshire between 19 november 2018 and 23 november 2018.'
    49 = 'The identification of initial need and risk in solihull.'
    50 = 'Serious youth violence in somerset.'
    51 = 'Abuse and neglect in south tyneside metropolitan borough between 22 and 26 february 2016.'
    52 = 'Child sexual exploitation, children associated with gangs and at risk of exploitation and children missing from home, care or education in southend-on-sea between 19 and 23 march 2018.'
    53 = 'Abuse and neglect in stockton-on-tees between 20 and 24 november 2017.'
    54 = 'Children and families who need help.'
    55 = 'Children and families in surrey who need help.'
    56 = 'Identification of initial need and risk in torbay.'
    57 = 'Identification of initial need and risk in walsall.'
    58 = 'Abuse and neglect in wiltshire between 31 october and 4 november 2016.'
    59 = 'Children and families who need help in the wirral.'
    60 = 'Abuse and neglect in wokingham borough council between 22 may and 26 may 2017.'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 11).Value = $updates[$row]
}
